$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Variables")

# Widen the "label" column on the Variables sheet and wrap its header text
$ws1.Columns.Item(3).ColumnWidth = 48.8
$ws1.Range("C1").WrapText = $true

# New data row describing the USbmi_m variable
$ws1.Range("B2").Value = "USbmi_m"
$ws1.Range("C2").Value = "Body-Mass-Index  Marker für die vergröberten Werte [kg/m2]"
$ws1.Range("D2").Value = "decimal"

# Wrap the long label text and grow the row to fit two lines
$ws1.Range("C2").WrapText = $true
$ws1.Rows.Item(2).RowHeight = 29

# Restore cursor position as it was left in the source file
$ws1.Range("C17").Select() | Out-Null

$wb.Save()
